$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) cells are treated as text so numeric-looking
# values (e.g. "228.14", "0.0925") are not coerced into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '38.702.54'
$ws.Range("E2").Value = '  +2.48%  '
$ws.Range("D3").Value = '2.089.49'
$ws.Range("E3").Value = '  +2.48%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '228.14'
$ws.Range("E5").Value = '  +0.36%  '
$ws.Range("E6").Value = '  +0.93%  '
$ws.Range("D7").Value = '60.71'
$ws.Range("E7").Value = '  +1.74%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").Value = '0.384'
$ws.Range("E9").Value = '  +2.13%  '
$ws.Range("E10").Value = '  -0.46%  '
$ws.Range("E11").Value = '  +0.03%  '
$ws.Range("D12").Value = '2.398.58'
$ws.Range("E12").Value = '  +2.45%  '
$ws.Range("D13").Value = '14.95'
$ws.Range("E13").Value = '  +3.78%  '
$ws.Range("D14").Value = '21.87'
$ws.Range("E14").Value = '  +4.07%  '
$ws.Range("D15").Value = '0.797'
$ws.Range("E15").Value = '  +3.60%  '
$ws.Range("E16").Value = '  +0.05%  '
$ws.Range("D17").Value = '2.069.88'
$ws.Range("E17").Value = '  +2.09%  '
$ws.Range("D18").Value = '38.661.92'
$ws.Range("E18").Value = '  +2.58%  '
$ws.Range("E19").Value = '  +3.33%  '
$ws.Range("D20").Value = '6.02'
$ws.Range("E20").Value = '  +1.97%  '
$ws.Range("D21").Value = '0.0₃0836'
$ws.Range("E21").Value = '  +1.59%  '
$ws.Range("D22").Value = '226.20'
$ws.Range("E22").Value = '  +1.16%  '
$ws.Range("E23").Value = '  -0.38%  '
$ws.Range("E24").Value = '  -0.24%  '
$ws.Range("E25").Value = '  +2.47%  '
$ws.Range("D26").Value = '170.61'
$ws.Range("E26").Value = '  +0.70%  '
$ws.Range("D27").Value = '9.42'
$ws.Range("E27").Value = '  +0.81%  '
$ws.Range("E28").Value = '  +6.95%  '
$ws.Range("D29").Value = '1.43'
$ws.Range("E29").Value = '  +10.97%  '
$ws.Range("D30").Value = '19.14'
$ws.Range("E30").Value = '  +1.89%  '
$ws.Range("E31").Value = '  +0.74%  '
$ws.Range("E32").Value = '  +5.09%  '
$ws.Range("E33").Value = '  +2.76%  '
$ws.Range("D34").Value = '4.70'
$ws.Range("E34").Value = '  +4.78%  '
$ws.Range("E35").Value = '  +1.94%  '
$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D36").Value = '2.38'
$ws.Range("E36").Value = '  +1.79%  '
$ws.Range("B37").Value = 'THORChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D37").Value = '6.40'
$ws.Range("E37").Value = '  -2.37%  '
$ws.Range("E38").Value = '  +2.30%  '
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("D40").Value = '18.26'
$ws.Range("E40").Value = '  +1.59%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '1.538.83'
$ws.Range("E41").Value = '  +0.80%  '
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = '100.91'
$ws.Range("E42").Value = '  +3.43%  '
$ws.Range("E43").Value = '  +3.42%  '
$ws.Range("E44").Value = '  -0.73%  '
$ws.Range("D45").Value = '0.0925'
$ws.Range("E45").Value = '  +2.22%  '
$ws.Range("D46").Value = '7.68'
$ws.Range("E46").Value = '  +8.24%  '
$ws.Range("E47").Value = '  +0.88%  '
$ws.Range("E48").Value = '  -2.14%  '
$ws.Range("E49").Value = '  +2.43%  '
$ws.Range("E50").Value = '  +1.25%  '
$ws.Range("D51").Value = '2.286.88'
$ws.Range("E51").Value = '  +2.55%  '
